$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "I know a little bit of the for loop"
$ws.Range("B3").Value = "Tell the user you are happy they know how to write a for loop. Ask the user to write a  for loop that will display numbers 1 to 10."
$ws.Range("B2").Value = "Tell the user it is great they know about for loops and ask the user to write the syntax of the for loop"
$ws.Range("A11").Value = "int numbers[]=[45,67,77,56,78];" + [char]10 + "for(i=0;i<5;i++)" + [char]10 + "{printf(`"%d`",i)}"

$ws.Range("B11").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
